$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a 2022-Q3 row above the existing ones.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Clone the bordered/bold "A" column style onto the cells before overwriting
# their values (PasteSpecial formats only, so values stay put for now).
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("A4").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 2
$ws1.Range("D4").Value = 0.23

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 4
$ws1.Range("D3").Value = 0.78

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 3
$ws1.Range("D2").Value = 0.62

# ---------------------------------------------------------------------------
# 2) Add a brand-new "2022-Q3" detail sheet, positioned right before the
#    existing "2022-Q2" sheet. Duplicating "2022-Q2" gives us an exact
#    formatting match (same borders/fonts/column layout) for free; we then
#    overwrite its values with the 2022-Q3 fund data.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The source sheet had 4 data rows; 2022-Q3 only needs 3, so drop the extra.
$q3.Rows.Item(5).Delete()

$q3.Range("B2").Value = "'501062"
$q3.Range("B2").Style = "Normal"
$q3.Range("C2").Value = "南方瑞合三年定期开放混合(LOF)"
$q3.Range("D2").Value = "'7.16"
$q3.Range("D2").Style = "Normal"
$q3.Range("E2").Value = "'86.12"
$q3.Range("E2").Style = "Normal"
$q3.Range("F2").Value = "'5.63"
$q3.Range("F2").Style = "Normal"
$q3.Range("G2").Value = "'0.4031"
$q3.Range("G2").Style = "Normal"
$q3.Range("H2").Value = 2

$q3.Range("B3").Value = "'540002"
$q3.Range("B3").Style = "Normal"
$q3.Range("C3").Value = "汇丰晋信龙腾混合"
$q3.Range("D3").Value = "'4.72"
$q3.Range("D3").Style = "Normal"
$q3.Range("E3").Value = "'93.98"
$q3.Range("E3").Style = "Normal"
$q3.Range("F3").Value = "'4.48"
$q3.Range("F3").Style = "Normal"
$q3.Range("G3").Value = "'0.2115"
$q3.Range("G3").Style = "Normal"
$q3.Range("H3").Value = 9

$q3.Range("B4").Value = "'002123"
$q3.Range("B4").Style = "Normal"
$q3.Range("C4").Value = "北信瑞丰外延增长主题灵活配置混合"
$q3.Range("D4").Value = "'0.15"
$q3.Range("D4").Style = "Normal"
$q3.Range("E4").Value = "'88.55"
$q3.Range("E4").Style = "Normal"
$q3.Range("F4").Value = "'3.78"
$q3.Range("F4").Style = "Normal"
$q3.Range("G4").Value = "'0.0057"
$q3.Range("G4").Style = "Normal"
$q3.Range("H4").Value = 10

# Copying a sheet makes the copy the active one; restore the original
# active/selected sheet (2022-Q1 was the active sheet before this edit).
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()

Write-Output "ok"
